# Update "paises.xlsx" COVID-19 stats worksheet with the newer snapshot
# (23 Abril 2020, 17:52) and refresh the "Casos totales" descending sort
# order, since Mali's updated count (309) now exceeds Isla de Man (307)
# and Venezuela (298).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp (row 1) -------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 23 de Abril de 2020 a las 17:52"

# --- Straightforward numeric refreshes (no re-sorting needed) ------------
# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 852610
$ws.Cells.Item(4,3).Value = 3893
$ws.Cells.Item(4,4).Value = 84191
$ws.Cells.Item(4,5).Value = 720124
$ws.Cells.Item(4,6).Value = 14344
$ws.Cells.Item(4,7).Value = 636
$ws.Cells.Item(4,8).Value = 48295

# Row 14: Brasil
$ws.Cells.Item(14,2).Value = 46701
$ws.Cells.Item(14,3).Value = 944
$ws.Cells.Item(14,5).Value = 18443
$ws.Cells.Item(14,7).Value = 34
$ws.Cells.Item(14,8).Value = 2940

# Row 34: Polonia
$ws.Cells.Item(34,2).Value = 10511
$ws.Cells.Item(34,3).Value = 342
$ws.Cells.Item(34,5).Value = 8317
$ws.Cells.Item(34,7).Value = 28
$ws.Cells.Item(34,8).Value = 454

# Row 44: Chequia
$ws.Cells.Item(44,2).Value = 7138
$ws.Cells.Item(44,3).Value = 6
$ws.Cells.Item(44,4).Value = 2152
$ws.Cells.Item(44,5).Value = 4776

# Row 58: Moldavia
$ws.Cells.Item(58,5).Value = 2185
$ws.Cells.Item(58,7).Value = 5
$ws.Cells.Item(58,8).Value = 80

# Row 81: Cuba
$ws.Cells.Item(81,2).Value = 1235
$ws.Cells.Item(81,3).Value = 46
$ws.Cells.Item(81,4).Value = 365
$ws.Cells.Item(81,5).Value = 827
$ws.Cells.Item(81,6).Value = 14
$ws.Cells.Item(81,7).Value = 3
$ws.Cells.Item(81,8).Value = 43

# Row 114: Sri Lanka
$ws.Cells.Item(114,2).Value = 340
$ws.Cells.Item(114,3).Value = 10
$ws.Cells.Item(114,5).Value = 226

# Row 115: Mauricio
$ws.Cells.Item(115,2).Value = 331
$ws.Cells.Item(115,3).Value = 2
$ws.Cells.Item(115,4).Value = 266
$ws.Cells.Item(115,5).Value = 56

# Row 153: Zambia
$ws.Cells.Item(153,2).Value = 76
$ws.Cells.Item(153,3).Value = 2
$ws.Cells.Item(153,4).Value = 37

# --- Re-sort block around Mali / Isla de Man / Venezuela -----------------
# Before: row119=Isla de Man(307), row120=Venezuela(298), row121=Mali(293)
# Mali's refreshed total (309) now outranks both, so it moves to row 119;
# Isla de Man and Venezuela (unchanged totals) shift down one row each.
$ws.Cells.Item(119,1).Value = "Mali"
$ws.Cells.Item(119,2).Value = 309
$ws.Cells.Item(119,3).Value = 16
$ws.Cells.Item(119,4).Value = 77
$ws.Cells.Item(119,5).Value = 211
$ws.Cells.Item(119,6).Value = 0
$ws.Cells.Item(119,7).Value = 4
$ws.Cells.Item(119,8).Value = 21

$ws.Cells.Item(120,1).Value = "Isla de Man"
$ws.Cells.Item(120,2).Value = 307
$ws.Cells.Item(120,3).Value = 0
$ws.Cells.Item(120,4).Value = 212
$ws.Cells.Item(120,5).Value = 80
$ws.Cells.Item(120,6).Value = 20
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 15

$ws.Cells.Item(121,1).Value = "Venezuela"
$ws.Cells.Item(121,2).Value = 298
$ws.Cells.Item(121,3).Value = 0
$ws.Cells.Item(121,4).Value = 122
$ws.Cells.Item(121,5).Value = 166
$ws.Cells.Item(121,6).Value = 4
$ws.Cells.Item(121,7).Value = 0
$ws.Cells.Item(121,8).Value = 10
